$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D16").Value = "Stitch it in Time: GAN-Based Facial Editing of Real Videos"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/189"

$ws.Range("D46").Value = "항암화학요법 부작용"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/457"
